$d = $word.ActiveDocument

# --- 1) Insert all new paragraphs (Entidades, Relaciones, Atributos, etc.) ---
#     Using Range.InsertXML with literal OOXML keeps the exact paragraph/run
#     layout (including the multi-run "Atributos" paragraph with its curly
#     quotes) instead of letting the host engine merge runs together.
$end = $d.Content
$end.Collapse(0)

$fragment = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Entidades:</w:t></w:r></w:p>
<w:p><w:r><w:t>Profesional-&gt;Tipo: El tipo es ‘interno’ o ‘externo’.</w:t></w:r></w:p>
<w:p><w:r><w:t>Especialidad: Representa una especialidad médica específica.</w:t></w:r></w:p>
<w:p><w:r><w:t>Lugar: Define la ubicación y dirección que tendrá el turno.</w:t></w:r></w:p>
<w:p><w:r><w:t>Procedimiento médico: Representa un procedimiento médico que se realiza sobre algún paciente. En él se establecen las condiciones en las que debe presentarse.</w:t></w:r></w:p>
<w:p><w:r><w:t>Block de Turnos: Representa un cronograma de turnos, es específico de cada médico y especialidad.</w:t></w:r></w:p>
<w:p><w:r><w:t>Relaciones:</w:t></w:r></w:p>
<w:p><w:r><w:t>Paciente-Necesita-Turno: La relación establece el paciente que es atendido por un médico.</w:t></w:r></w:p>
<w:p><w:r><w:t>Profesional-Solicita-Turno de quirófano y cama: La relación establece el médico que reserva un turno de cama o quirófano que corresponde a un paciente.</w:t></w:r></w:p>
<w:p><w:r><w:t>Cobertura-Cubre-Procedimiento Médico: Representa la cobertura de un plan de obra social a un procedimiento médico.</w:t></w:r></w:p>
<w:p><w:r><w:t>Turno-Turno en Block-Block de Turnos: Representa el requerimiento de encasillar a cada turno en un block de turnos.</w:t></w:r></w:p>
<w:p><w:r><w:t>Atributos:</w:t></w:r></w:p>
<w:p><w:r><w:t>Turno de servicio de diagnostico-&gt;Tipo de turno: El tipo de turno puede ser ‘</w:t></w:r><w:r><w:t>primera vez</w:t></w:r><w:r><w:t>’</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>‘</w:t></w:r><w:r><w:t>segunda vez</w:t></w:r><w:r><w:t>’</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>o</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>‘</w:t></w:r><w:r><w:t>demandas espontáneas</w:t></w:r><w:r><w:t>’.</w:t></w:r></w:p>
<w:p><w:r><w:t>Block de Turnos-&gt;Tipo de Agenda: Estable que cual es el tipo de atención, es decir, si es personalizada o en grupos de N personas.</w:t></w:r></w:p>
<w:p><w:r><w:t>Block de Turnos-&gt;Bloqueado: Este atributo indica si se pueden agregar o no turnos a este block.</w:t></w:r></w:p>
<w:p><w:r><w:t>Turno-&gt;Monto a abonar: ¿?</w:t></w:r></w:p>
<w:p><w:r><w:t>Cubre-&gt;</w:t></w:r><w:r><w:t>Exención: Es el porcentaje que se descuenta del precio de lista del procedimiento médico.</w:t></w:r></w:p>
<w:p><w:r><w:t>Cubre-&gt;</w:t></w:r><w:r><w:t xml:space="preserve">Valor Copago: Es el monto que tiene que pagar el paciente. Exención y Valor Copago </w:t></w:r><w:r><w:t>representan lo mismo, sólo uno toma valor.</w:t></w:r></w:p>
<w:p><w:r><w:t>Cubre-&gt;</w:t></w:r><w:r><w:t>Documentación:</w:t></w:r><w:r><w:t xml:space="preserve"> Tiene los documentos que necesita presentar el paciente para lograr la cobertura.</w:t></w:r></w:p>
'@
$end.InsertXML($fragment)

# --- 2) Leave behind a numbering definition (word/numbering.xml) and the
#     "List Paragraph" style, mirroring the artifacts Word creates the first
#     time a bulleted list is tried out, even though no paragraph above ends
#     up using bullets. ---
$scratchAnchor = $d.Content
$scratchAnchor.Collapse(0)
$scratchAnchor.InsertParagraphAfter()
$scratchPara = $d.Paragraphs.Last
$scratchPara.Range.Text = "scratch"
$scratchPara.Range.ListFormat.ApplyBulletDefault()
$scratchPara.Range.Delete()

$listStyle = $d.Styles.Add("Prrafodelista", 1)
$listStyle.NameLocal = "List Paragraph"
$listStyle.BaseStyle = $d.Styles.Item("Normal")
$listStyle.Priority = 34
$listStyle.QuickStyle = $true
$listStyle.ParagraphFormat.LeftIndent = 36
$listStyle.NoSpaceBetweenParagraphsOfSameStyle = $true

Write-Host ("Paragraph count: " + $d.Paragraphs.Count)
